$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.885.87"
$ws.Range("E2").Value = "  -0.51%  "
Set-TextValue "D3" "2.567.80"
$ws.Range("E3").Value = "  -0.14%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.10%  "
Set-TextValue "D5" "584.78"
$ws.Range("E5").Value = "  +0.03%  "
Set-TextValue "D6" "143.73"
$ws.Range("E6").Value = "  -2.48%  "
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("E9").Value = "  -2.68%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  -2.42%  "
Set-TextValue "D13" "26.95"
$ws.Range("E13").Value = "  -1.92%  "
Set-TextValue "D14" "3.025.39"
$ws.Range("E14").Value = "  -0.32%  "
Set-TextValue "D15" "62.774.77"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("E16").Value = "  -2.20%  "
Set-TextValue "D17" "2.584.84"
$ws.Range("E17").Value = "  +0.09%  "
Set-TextValue "D18" "11.02"
$ws.Range("E18").Value = "  -2.95%  "
Set-TextValue "D19" "339.65"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("E22").Value = "  +0.19%  "
Set-TextValue "D23" "5.72"
Set-TextValue "D24" "67.49"
$ws.Range("E24").Value = "  +0.98%  "
Set-TextValue "D25" "1.59"
$ws.Range("E25").Value = "  +7.17%  "
Set-TextValue "D26" "1.59"
$ws.Range("E26").Value = "  -2.96%  "
Set-TextValue "D27" "0.164"
$ws.Range("E27").Value = "  -3.96%  "
Set-TextValue "D28" "7.97"
$ws.Range("E28").Value = "  -2.31%  "
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("E30").Value = "  -3.20%  "
Set-TextValue "D31" "1.92"
$ws.Range("E31").Value = "  -2.85%  "
Set-TextValue "D32" "459.85"
$ws.Range("E32").Value = "  -0.51%  "
Set-TextValue "D33" "0.0₃0794"
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("E34").Value = "  +1.03%  "
Set-TextValue "D35" "176.43"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("E36").Value = "  +0.02%  "
Set-TextValue "D37" "0.396"
$ws.Range("E37").Value = "  -2.72%  "
Set-TextValue "D38" "18.79"
$ws.Range("E38").Value = "  -2.29%  "
Set-TextValue "D39" "4.52"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E41").Value = "  -3.47%  "
Set-TextValue "D42" "39.94"
$ws.Range("E42").Value = "  +0.96%  "
Set-TextValue "D43" "157.26"
$ws.Range("E43").Value = "  +3.99%  "
Set-TextValue "D44" "3.67"
$ws.Range("E44").Value = "  -3.56%  "
$ws.Range("E45").Value = "  +0.50%  "
Set-TextValue "D46" "0.629"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("E48").Value = "  -2.37%  "
Set-TextValue "D49" "0.0235"
$ws.Range("E49").Value = "  -1.71%  "
Set-TextValue "D50" "17.97"
$ws.Range("E50").Value = "  -2.87%  "
Set-TextValue "D51" "11.39"
$ws.Range("E51").Value = "  +0.07%  "
